$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Trial": update row 3 values/formats and append rows 4 and 5
# ---------------------------------------------------------------------
$trial = $wb.Worksheets.Item("Trial")

# --- Row 3 updates -----------------------------------------------------
# D3 and G3 switch from the "General" style (s=3) to the "0.00" style (s=4)
$trial.Range("D3").NumberFormat = "0.00"
$trial.Range("D3").Value = 0.4
$trial.Range("E3").Value = 0.77
$trial.Range("G3").NumberFormat = "0.00"
$trial.Range("G3").Value = 0.4
$trial.Range("H3").Value = 0.77
# I3 switches from "General" style (s=3) to the "0.00" style (s=4)
$trial.Range("I3").NumberFormat = "0.00"
$trial.Range("I3").Value = 0.6
$trial.Range("J3").Value = 0.001493
$trial.Range("L3").Value = 0.001789

# --- Row 4 (new) ---------------------------------------------------------
$trial.Rows.Item(4).RowHeight = 35.05

$trial.Range("A4").NumberFormat = "General"
$trial.Range("A4").Value = 1
$trial.Range("B4").NumberFormat = "General"
$trial.Range("B4").Value = 2
$trial.Range("C4").NumberFormat = "General"
$trial.Range("C4").Value = "RFC"
$trial.Range("D4").NumberFormat = "0.00"
$trial.Range("D4").Value = 0.45
$trial.Range("E4").NumberFormat = "0.00"
$trial.Range("E4").Value = 0.8
$trial.Range("F4").NumberFormat = "General"
$trial.Range("F4").Value = "over-fitting (high variance)"
$trial.Range("G4").NumberFormat = "0.00"
$trial.Range("G4").Value = 0.45
$trial.Range("H4").NumberFormat = "0.00"
$trial.Range("H4").Value = 0.8
$trial.Range("I4").NumberFormat = "0.00"
$trial.Range("I4").Value = 0.55
$trial.Range("J4").NumberFormat = "0.000000"
$trial.Range("J4").Value = 0.001342
$trial.Range("K4").NumberFormat = "0.00"
$trial.Range("K4").Value = 0.9
$trial.Range("L4").NumberFormat = "0.000000"
$trial.Range("L4").Value = 0.000575

# --- Row 5 (new) ---------------------------------------------------------
$trial.Rows.Item(5).RowHeight = 35.05

$trial.Range("A5").NumberFormat = "General"
$trial.Range("A5").Value = 1
$trial.Range("B5").NumberFormat = "General"
$trial.Range("B5").Value = 3
$trial.Range("C5").NumberFormat = "General"
$trial.Range("C5").Value = "RFC"
$trial.Range("D5").NumberFormat = "0.00"
$trial.Range("D5").Value = 0.45
$trial.Range("E5").NumberFormat = "0.00"
$trial.Range("E5").Value = 0.8
$trial.Range("F5").NumberFormat = "General"
$trial.Range("F5").Value = "over-fitting (high variance)"
$trial.Range("G5").NumberFormat = "0.00"
$trial.Range("G5").Value = 0.45
$trial.Range("H5").NumberFormat = "0.00"
$trial.Range("H5").Value = 0.8
$trial.Range("I5").NumberFormat = "0.00"
$trial.Range("I5").Value = 0.55
$trial.Range("J5").NumberFormat = "0.000000"
$trial.Range("J5").Value = 0.001138
$trial.Range("K5").NumberFormat = "0.00"
$trial.Range("K5").Value = 0.9
$trial.Range("L5").NumberFormat = "0.000000"
$trial.Range("L5").Value = 0.000831

# Update the remembered selection for this sheet (was A3, now J4)
$trial.Range("J4").Select()

# ---------------------------------------------------------------------
# Sheet "RFC_params": rewrite the boolean cells as TRUE()/FALSE() formulas
# and move the remembered selection
# ---------------------------------------------------------------------
$rfc = $wb.Worksheets.Item("RFC_params")

$rfc.Range("B3").Formula = "=TRUE()"
$rfc.Range("O3").Formula = "=FALSE()"

$rfc.Range("B3").Select()
